$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

# Write the discount unit price value to the "Client discount" row
$ws.Range("E18").Value = 120

# Write the signature (name, email address) to the cell
$ws.Range("A31").Value = "RPA Dev, developer.rpa@mail.com"
